$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "24/10/2025"
$ws.Range("B12").Value = "Gorica"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "Hajduk Split"
$ws.Range("F12").Value = "L"
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 1.56
$ws.Range("L12").Value = 2.2
$ws.Range("M12").Value = 9
$ws.Range("N12").Value = 13
$ws.Range("O12").Value = 3
$ws.Range("P12").Value = 8
